$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# ------------------------------------------------------------------
# Row 3: remove the two stray formatted-but-empty cells F3 and G3
# ------------------------------------------------------------------
$ws.Range("F3:G3").Clear()

# ------------------------------------------------------------------
# Row 20: turn into a merged "Mensajes" banner (style = Incorrecto + centered)
# ------------------------------------------------------------------
$ws.Range("A20:J20").Style = "Incorrecto"
$ws.Range("A20:J20").HorizontalAlignment = -4108
$ws.Range("A20").Value = "Mensajes"
$ws.Range("A20:J20").Merge()

# ------------------------------------------------------------------
# Row 32: turn into a merged "dibujos" banner (style = Incorrecto + centered)
# ------------------------------------------------------------------
$ws.Range("A32:I32").Style = "Incorrecto"
$ws.Range("A32:I32").HorizontalAlignment = -4108
$ws.Range("A32").Value = "dibujos"
$ws.Range("A32:I32").Merge()

# ------------------------------------------------------------------
# Row 33: clear the old "Tag" value out of G33 (keep its centered style),
# add the new "Tags" header in J33
# ------------------------------------------------------------------
$ws.Range("G33").ClearContents()
$ws.Range("J33").Value = "Tags"

# ------------------------------------------------------------------
# Row 34: rework the small "imagen" key table and start a new tags table
# ------------------------------------------------------------------
$ws.Range("B34").Value = "ID-galeria"
$ws.Range("C34").Value = "ID-usuario"

$ws.Range("D34").Value = "ID-box-tags-imagen"
$ws.Range("E34").Value = ""
$ws.Range("F34").Value = "ID-box-tags"
$ws.Range("D34:F34").HorizontalAlignment = 1

$ws.Range("G34").ClearContents()
$ws.Range("H34:I34").Clear()

$ws.Range("J34").Value = "ID-tags"
$ws.Range("J34").HorizontalAlignment = -4108
$ws.Range("K34").Value = "nombre-tag"

# ------------------------------------------------------------------
# Rows 35-39: "dibujos" rows (A/D/F) and the new "tags" lookup table (J/K)
# (the old rows 35/36 only held a lone formatted G cell - drop those first)
# ------------------------------------------------------------------
$ws.Range("G35:G36").Clear()

$ws.Range("A35").Value = 1
$ws.Range("D35").Value = 1
$ws.Range("F35").Value = 1
$ws.Range("J35").Value = 1
$ws.Range("J35").HorizontalAlignment = -4108
$ws.Range("K35").Value = "dibujo"

$ws.Range("A36").Value = 2
$ws.Range("J36").Value = 2
$ws.Range("J36").HorizontalAlignment = -4108
$ws.Range("K36").Value = "arte conceptual"

$ws.Range("A37").Value = 3
$ws.Range("J37").Value = 3
$ws.Range("J37").HorizontalAlignment = -4108
$ws.Range("K37").Value = "renacentista"

$ws.Range("A38").Value = 4
$ws.Range("J38").Value = 4
$ws.Range("J38").HorizontalAlignment = -4108
$ws.Range("K38").Value = "diseño de entorno"

$ws.Range("J39").Value = 5
$ws.Range("J39").HorizontalAlignment = -4108
$ws.Range("K39").Value = "videojuegos"

# ------------------------------------------------------------------
# Row 43: two empty centered cells
# ------------------------------------------------------------------
$ws.Range("A43:B43").Style = "Normal"
$ws.Range("A43:B43").HorizontalAlignment = -4108

# ------------------------------------------------------------------
# View: drop the old scrolled-down selection, select H29 instead
# ------------------------------------------------------------------
$ws.Activate()
$ws.Range("H29").Select()
